{"js": "// Replace each \"dividend\u00f7divisor=\" exercise text in the table with its\n// new value. Every old string in this document is unique, so a scoped\n// Find & Replace (matchCase, whole match via the trailing \"=\") is safe.\nconst replacements = [\n  [\"514\u00f78=\", \"869\u00f75=\"],\n  [\"277\u00f76=\", \"562\u00f79=\"],\n  [\"299\u00f72=\", \"944\u00f72=\"],\n  [\"522\u00f77=\", \"377\u00f74=\"],\n  [\"853\u00f75=\", \"303\u00f77=\"],\n  [\"786\u00f72=\", \"653\u00f75=\"],\n  [\"355\u00f75=\", \"823\u00f79=\"],\n  [\"200\u00f77=\", \"577\u00f75=\"],\n  [\"889\u00f75=\", \"195\u00f77=\"],\n  [\"285\u00f79=\", \"437\u00f74=\"],\n  [\"755\u00f72=\", \"195\u00f78=\"],\n  [\"647\u00f79=\", \"244\u00f74=\"],\n  [\"546\u00f78=\", \"861\u00f73=\"],\n  [\"476\u00f75=\", \"455\u00f75=\"],\n  [\"243\u00f77=\", \"273\u00f72=\"],\n  [\"736\u00f76=\", \"679\u00f73=\"],\n  [\"437\u00f76=\", \"315\u00f78=\"],\n  [\"661\u00f73=\", \"757\u00f73=\"],\n  [\"486\u00f78=\", \"991\u00f79=\"],\n  [\"138\u00f74=\", \"724\u00f73=\"],\n  [\"388\u00f78=\", \"694\u00f78=\"],\n  [\"648\u00f77=\", \"414\u00f75=\"],\n  [\"850\u00f78=\", \"662\u00f78=\"],\n  [\"329\u00f79=\", \"292\u00f73=\"],\n  [\"677\u00f72=\", \"524\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"dividend\u00f7divisor=\" exercise text in the table with its\n# new value, via Find & Replace over the whole document content.\n# Every \"old\" string in this document is unique, so a single\n# Find.Execute(..., Replace:=wdReplaceAll) per pair is safe and will\n# touch exactly one run.\n\n$replacements = @(\n  @(\"514\u00f78=\", \"869\u00f75=\"),\n  @(\"277\u00f76=\", \"562\u00f79=\"),\n  @(\"299\u00f72=\", \"944\u00f72=\"),\n  @(\"522\u00f77=\", \"377\u00f74=\"),\n  @(\"853\u00f75=\", \"303\u00f77=\"),\n  @(\"786\u00f72=\", \"653\u00f75=\"),\n  @(\"355\u00f75=\", \"823\u00f79=\"),\n  @(\"200\u00f77=\", \"577\u00f75=\"),\n  @(\"889\u00f75=\", \"195\u00f77=\"),\n  @(\"285\u00f79=\", \"437\u00f74=\"),\n  @(\"755\u00f72=\", \"195\u00f78=\"),\n  @(\"647\u00f79=\", \"244\u00f74=\"),\n  @(\"546\u00f78=\", \"861\u00f73=\"),\n  @(\"476\u00f75=\", \"455\u00f75=\"),\n  @(\"243\u00f77=\", \"273\u00f72=\"),\n  @(\"736\u00f76=\", \"679\u00f73=\"),\n  @(\"437\u00f76=\", \"315\u00f78=\"),\n  @(\"661\u00f73=\", \"757\u00f73=\"),\n  @(\"486\u00f78=\", \"991\u00f79=\"),\n  @(\"138\u00f74=\", \"724\u00f73=\"),\n  @(\"388\u00f78=\", \"694\u00f78=\"),\n  @(\"648\u00f77=\", \"414\u00f75=\"),\n  @(\"850\u00f78=\", \"662\u00f78=\"),\n  @(\"329\u00f79=\", \"292\u00f73=\"),\n  @(\"677\u00f72=\", \"524\u00f73=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $oldText\"\n  }\n}\n"}
